$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# Update R2 resistor value: 165k (0,1%) -> 110k (0,1%), clear its LCSC part number (E20)
$ws.Range("D20").Value = "110k (0,1%)"
$ws.Range("E20").ClearContents()

# Update R4 resistor value: 15k (0.1%) -> 10k (0.1%), clear its LCSC part number (E21)
$ws.Range("D21").Value = "10k (0.1%)"
$ws.Range("E21").ClearContents()

# Adjust the window / view state to match the saved workbook view
$ws.Activate()
$ws.Range("H21").Select()
